$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A26").Value = "V3 MP Firmware"
$ws.Range("B26").Value = "7.1.1"
$ws.Range("C26").Value = (Get-Date -Year 2022 -Month 7 -Day 18).Date
$ws.Range("D26").Value = "Zound_Hendrix_M_Lite_V3_hwEVT_btswv7.1.0_20220718"
$ws.Range("E26").Value = "7.1.1"
$ws.Range("G26").Value = 0.6
$ws.Range("H26").Value = 3.1
$ws.Range("J26").Value = "N/A"
$ws.Range("K26").Value = "Modify the amp RT9120S output PWM to 768KHz to lower the power consumption."

$ws.Rows.Item(26).RowHeight = 32.799999999999997

$ws.Range("K27").Select()
